$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.256.70'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.268.30'
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.94'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.32'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.529'
$ws.Range('E7').Value = '  -0.89%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.20'
$ws.Range('E10').Value = '  -1.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('E11').Value = '  -2.57%  '
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.96'
$ws.Range('E13').Value = '  +2.90%  '
$ws.Range('D14').Value = '2.619.92'
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('D16').Value = '2.280.87'
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('D18').Value = '42.116.70'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.31'
$ws.Range('E19').Value = '  -4.37%  '
$ws.Range('E20').Value = '  -1.70%  '
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.93'
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.01'
$ws.Range('E23').Value = '  -2.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.58'
$ws.Range('E24').Value = '  -2.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.97'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.58'
$ws.Range('E27').Value = '  -2.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.83'
$ws.Range('E28').Value = '  +2.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.56'
$ws.Range('E29').Value = '  -1.90%  '
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '162.57'
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('E32').Value = '  -2.78%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.71'
$ws.Range('E35').Value = '  +0.92%  '
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('E38').Value = '  -4.49%  '
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('E41').Value = '  -3.45%  '
$ws.Range('E42').Value = '  +3.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.42'
$ws.Range('D44').Value = '1.951.77'
$ws.Range('E44').Value = '  -3.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0282'
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.93'
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.86'
$ws.Range('E47').Value = '  -4.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.69'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '92.57'
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('E50').Value = '  -1.86%  '
$ws.Range('E51').Value = '  -3.14%  '
